# Data analysis for Biology completed
# Adds a "total" summary (Male/Female sums across the two replicate blocks)
# to the Vial2 and Vial4 sheets, then leaves the view on the Vial5 tab.

$wb = $excel.ActiveWorkbook

# ---- Vial2 (sheet1) ----------------------------------------------------
$ws1 = $wb.Worksheets.Item("Vial2")

$ws1.Range("F3").Value = "G2"
$ws1.Range("G3").Value = "Male"
$ws1.Range("H3").Value = "Female"

$ws1.Range("F4").Value = "Red"
$ws1.Range("G4").Formula = "=SUM(B4,B8)"
$ws1.Range("H4").Formula = "=SUM(C4,C8)"

$ws1.Range("F5").Value = "White"
$ws1.Range("G5").Formula = "=SUM(B5,B9)"
$ws1.Range("H5").Formula = "=SUM(C5,C9)"

[void]$ws1.Range("G4:H5").Select()

# ---- Vial3 (sheet2) -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Vial3")
[void]$ws2.Range("B3:C4").Select()

# ---- Vial4 (sheet3) ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Vial4")

$ws3.Range("H2").Value = "Male"
$ws3.Range("I2").Value = "Female"

$ws3.Range("G3").Value = "Red"
$ws3.Range("H3").Formula = "=SUM(B3,B6)"
$ws3.Range("I3").Formula = "=SUM(C3,C6)"

$ws3.Range("G4").Value = "White"
$ws3.Range("H4").Formula = "=SUM(B4,B7)"
$ws3.Range("I4").Formula = "=SUM(C4,C7)"

[void]$ws3.Range("H3:I4").Select()

# ---- Vial5 (sheet4) -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Vial5")
[void]$ws4.Range("D16").Select()
[void]$ws4.Activate()
